# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet: Home row (row 2)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 245
$wsOff.Range("C2").Value = 178
$wsOff.Range("D2").Value = 63
$wsOff.Range("E2").Value = 20
$wsOff.Range("F2").Value = 5

# DEF sheet: Home row (row 2)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 364
$wsDef.Range("C2").Value = 245
$wsDef.Range("D2").Value = 72
$wsDef.Range("E2").Value = 28
$wsDef.Range("F2").Value = 7
$wsDef.Range("G2").Value = 3
